# Applies the "Finish code restructuring for employees" edit to the
# Product features workbook:
#   - product backlog: row 12 and row 14 marked as "(3) Completed"
#     (highlighted like the other finished rows, using the same
#     yellow-fill/word-wrap formatting already used on rows 2,4,6,7,9,10)
#   - product backlog: three new backlog items appended (rows 17-19)
#   - product backlog: row 19 grows taller to fit its wrapped text
#   - selection/cursor position updated on "product backlog" and
#     "sprint backlog" sheets

$wb = $excel.ActiveWorkbook

$products = $wb.Worksheets.Item("product backlog")
$sprint   = $wb.Worksheets.Item("sprint backlog")

$yellow = 65535
$completed = "(3) Completed"

# --- Row 12: "setup virtual environment for the app" is now finished ---
$row12 = $products.Range("A12:I12")
$row12.Interior.Color = $yellow
$row12.WrapText = $true
$products.Range("I12").Value = $completed

# --- Row 14: "improve code design usign design patterns" is now finished ---
$row14 = $products.Range("A14:I14")
$row14.Interior.Color = $yellow
$row14.WrapText = $true
$products.Range("I14").Value = $completed

# --- New backlog items ---
$products.Range("E17").Value = "add photo for users"
$products.Range("E18").Value = "add logo for organizations"
$products.Range("E19").Value = "delete View button on user form. Rename " + [char]0x0395 + [char]0x03C0 + [char]0x03B5 + [char]0x03BE + [char]0x03B5 + [char]0x03C1 + [char]0x03B3 + [char]0x03B1 + [char]0x03C3 + [char]0x03AF + [char]0x03B1 + " to " + [char]0x03A0 + [char]0x03C1 + [char]0x03BF + [char]0x03B2 + [char]0x03BF + [char]0x03BB + [char]0x03AE + "/" + [char]0x0395 + [char]0x03C0 + [char]0x03B5 + [char]0x03BE + [char]0x03B5 + [char]0x03C1 + [char]0x03B3 + [char]0x03B1 + [char]0x03C3 + [char]0x03AF + [char]0x03B1

# Row 19 wraps onto three lines, so it needs to be taller.
$products.Rows.Item(19).RowHeight = 45

# --- Update remembered selections ---
$products.Range("E20").Select()

$sprint.Activate()
$sprint.Range("E2").Select()
$products.Activate()
